$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ur = $ws.UsedRange
    foreach ($cell in $ur.Cells) {
        $v = $cell.Value2
        if ($v -is [string] -and $v.Contains(",")) {
            $newVal = $v -replace ',(?!\s)', ', '
            if ($newVal -ne $v) {
                $cell.Value2 = $newVal
            }
        }
    }
}
